$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.815.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.249.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.99"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.640"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +10.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.98"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.55"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0958"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.32"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.574.88"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.896"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.79"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.246.77"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.585.93"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0977"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.97"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.98"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.62"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.29%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.62%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.95%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.52"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.90"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.10"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +11.98%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0789"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.124"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.77"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.68"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.16"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.27"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.56"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.77"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.65"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.02"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.46%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.39"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.37%  "

